$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the year/value table by one column (L), mirroring column K's
# formatting, and add the new 2021 / 269 data point.

# L3: blank bottom-border cell matching K3's style
$ws.Range("K3").Copy() | Out-Null
$ws.Range("L3").PasteSpecial(-4122) | Out-Null

# L4: year header cell matching K4's style, value 2021
$ws.Range("K4").Copy() | Out-Null
$ws.Range("L4").PasteSpecial(-4122) | Out-Null
$ws.Range("L4").Value = 2021

# L5: data cell matching K5's style, value 269
$ws.Range("K5").Copy() | Out-Null
$ws.Range("L5").PasteSpecial(-4122) | Out-Null
$ws.Range("L5").Value = 269

$excel.CutCopyMode = 0

# Update the selection to reflect the new active cell noted in the diff.
$ws.Range("N3").Select() | Out-Null
